$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 20: "user-page.php / searchRoom.php ..." fix entry is now resolved
#   - B20 switches from a real date value to the date written as plain text
#   - D20 gets a small wording fix (space added after "1.")
# ---------------------------------------------------------------------
$ws.Range("A20").Copy()
$ws.Range("B20").PasteSpecial(-4122)   # xlPasteFormats - reuse the plain "General" style
$ws.Range("B20").Value = "30.11.2015."
$ws.Range("D20").Value = "1. Jāatveras pareizajām profila lapām"

# ---------------------------------------------------------------------
# Row 21 used to be an empty template row (only the running number 20
# was filled in) - it now holds a real entry.
# ---------------------------------------------------------------------
$ws.Range("B21").Value = "02.12.2015."
$ws.Range("C21").Value = "searchRoom.php`nsearchCourse.php`nsearchPerson.php"
$ws.Range("D21").Value = "1. Ja lietotājs nav pievienojis porfila bildi, nepieciešams attēlot noklusējuma bildi"
$ws.Range("E21").Value = "J"
$ws.Range("C21").HorizontalAlignment = -4131   # xlLeft, matches D21's existing look
$ws.Rows("21").RowHeight = 45

# ---------------------------------------------------------------------
# Insert a brand-new row 22 (pushes the old blank/footer rows down by
# one) and fill it in with the next entry.
# ---------------------------------------------------------------------
$ws.Rows("22").Insert()

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "02.12.2015."
$ws.Range("C22").Value = "user-page.php`nsearchRoom.php`nsearchCourse.php`nsearchPerson.php"
$ws.Range("D22").Value = '1. Lauku "Lietotāja loma" vajag virs laukiem "apgūtie kursi, iegūtie diplomi, iegūtie sertifikāti, pasniedzamie kursi"'
$ws.Range("E22").Value = "J"
$ws.Range("C22").HorizontalAlignment = -4131   # xlLeft, matches D22's look
$ws.Rows("22").RowHeight = 60

# ---------------------------------------------------------------------
# Keep the view roughly where the author left it.
# ---------------------------------------------------------------------
$ws.Range("F22").Select()
